$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Worksheet")

# Add the new translated name/localization cells introduced alongside the
# existing character name ("Sina" next to "シィナ") and the rewritten
# line label ("Rewrite     -   Lily" next to "リリー").
$ws.Range("C3").Value = "Sina"
$ws.Range("C4").Value = "Rewrite     -   Lily"
